# Apply updated "dSF" (column F) values to the data rows that changed
# when the source data was re-pulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = -5
    15 = 1
    19 = -4
    20 = -1
    24 = -1
    25 = 3
    33 = -3
    56 = -1
    61 = -3
    62 = -2
    66 = 7
    67 = -10
    70 = -5
    73 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
